# Applies the "Trade #62 closed ..." update to live_trading_results.xlsx
#
# Summary of the change (per the target diff):
#  - momentum trade #40 (row 12 on the "momentum" sheet) transitions from
#    OPEN -> CLOSED with exit price / P&L / exit reason / duration filled in.
#  - Because trade #40 is now closed, it is appended as a new row on the
#    "All Trades" sheet (row 41).
#  - A brand new leadlag trade #62 is appended (OPEN) as row 51 on the
#    "leadlag" sheet.
#  - The "Summary" and "Comparison" aggregate sheets are recalculated to
#    reflect the newly closed trade.

function Set-TextValue {
    # Forces a value to be written as literal text, bypassing Excel's
    # automatic number / date / time / percentage detection so that values
    # such as "65.0%", "3.70", "2026-02-16" or "21:34:14" are stored exactly
    # as given instead of being silently coerced into numbers/dates.
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "momentum" sheet - close out trade #40 (row 12)
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

$momentum.Range("G12").Value = 67858.09089399999
Set-TextValue $momentum.Range("H12") "CLOSED"
$momentum.Range("I12").Value = 1.1631
$momentum.Range("J12").Value = 11.63
Set-TextValue $momentum.Range("M12") "time_exit_5min"
$momentum.Range("N12").Value = 5

# ---------------------------------------------------------------------
# 2) "All Trades" sheet - append the now-closed trade #40 as row 41
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("A41").Value = 40
Set-TextValue $allTrades.Range("B41") "2026-02-16"
Set-TextValue $allTrades.Range("C41") "21:29:12"
Set-TextValue $allTrades.Range("D41") "momentum"
Set-TextValue $allTrades.Range("E41") "DOWN"
$allTrades.Range("F41").Value = 68656.63
$allTrades.Range("G41").Value = 67858.09089399999
Set-TextValue $allTrades.Range("H41") "CLOSED"
$allTrades.Range("I41").Value = 1.1631
$allTrades.Range("J41").Value = 11.63
$allTrades.Range("K41").Value = 0.9
Set-TextValue $allTrades.Range("L41") "Downward momentum: -0.409% over 10 samples"
Set-TextValue $allTrades.Range("M41") "time_exit_5min"
$allTrades.Range("N41").Value = 5

# ---------------------------------------------------------------------
# 3) "leadlag" sheet - append new OPEN trade #62 as row 51
# ---------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

$leadlag.Range("A51").Value = 62
Set-TextValue $leadlag.Range("B51") "2026-02-16"
Set-TextValue $leadlag.Range("C51") "21:34:14"
Set-TextValue $leadlag.Range("D51") "leadlag"
Set-TextValue $leadlag.Range("E51") "DOWN"
$leadlag.Range("F51").Value = 68719.61500000001
# G51 stays blank (exit price not known yet - trade is still open)
Set-TextValue $leadlag.Range("H51") "OPEN"
$leadlag.Range("I51").Value = 0
$leadlag.Range("J51").Value = 0
$leadlag.Range("K51").Value = 0.7332
Set-TextValue $leadlag.Range("L51") "Coinbase leading with -0.073% move"
# M51 stays blank (no exit reason yet - trade is still open)
$leadlag.Range("N51").Value = 0

# ---------------------------------------------------------------------
# 4) "Summary" sheet - recompute aggregate stats
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

# Row 2: OVERALL / ALL COMBINED
$summary.Range("C2").Value = 40
Set-TextValue $summary.Range("D2") "65.0%"
Set-TextValue $summary.Range("E2") "+8.8308%"
Set-TextValue $summary.Range("F2") "+0.2208%"

# Row 4: STRATEGY / momentum
Set-TextValue $summary.Range("D4") "75.0%"
Set-TextValue $summary.Range("E4") "+4.2016%"
Set-TextValue $summary.Range("F4") "+0.3501%"

# ---------------------------------------------------------------------
# 5) "Comparison" sheet - recompute momentum strategy stats
# ---------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

Set-TextValue $comparison.Range("C3") "75.0%"
Set-TextValue $comparison.Range("D3") "4.74"
Set-TextValue $comparison.Range("E3") "+0.5918%"
Set-TextValue $comparison.Range("G3") "1.05"
